$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value.
# Update every data row (2 through 36) from 46061 (2026-02-08) to 46062 (2026-02-09).
for ($row = 2; $row -le 36; $row++) {
    $ws.Cells.Item($row, 3).Value = 46062
}
